$wb = $excel.ActiveWorkbook

# --- Add the new "Data" worksheet as the last (3rd) tab ---------------------
$gsmArena = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Data"

# --- Cell values --------------------------------------------------------
# Fill column C first (arg1, arg2, arg3) then the header (arg) last so the
# new shared-string table entries land in the same order the source file
# uses: arg1, arg2, arg3, arg.
$ws.Range("C2").Value = "arg1"
$ws.Range("C3").Value = "arg2"
$ws.Range("C4").Value = "arg3"
$ws.Range("C1").Value = "arg"

$ws.Range("A1").Value = "Execute"
$ws.Range("B1").Value = "TUID"

$ws.Range("A2").Value = "Y"
$ws.Range("A3").Value = "Y"
$ws.Range("A4").Value = "Y"

# --- Formatting: reuse the look already used by the other sheets ------------
# Bold, centered header style (matches the existing "Execute"/"TUID" headers).
$gsmArena.Range("A1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("C1").PasteSpecial(-4122)

# Centered "text" style used for the rest of the data cells.
$gsmArena.Range("A2").Copy() | Out-Null
$ws.Range("A2:A4").PasteSpecial(-4122)
$ws.Range("C2:C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# TUID header - bold + centered, General number format.
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").NumberFormat = "General"
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4108

# TUID values are numbers-as-text ("1", "2", "3"), centered, General format.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "3"
$ws.Range("B2:B4").NumberFormat = "General"
$ws.Range("B2:B4").HorizontalAlignment = -4108
$ws.Range("B2:B4").VerticalAlignment = -4108

# --- Column widths / sheet look -----------------------------------------
$ws.Columns.Item(2).ColumnWidth = 15.28515625
$ws.Columns.Item(3).ColumnWidth = 10.42578125

# --- Selection / activation ------------------------------------------------
$ws.Range("C2").Select()
$ws.Activate()
$gsmArena.Select()
$ws.Select()
